# Natmi following Dr Hou advice
# Replace the 3-cluster (FAPs/sCs) x 1-cluster (ECs) LR-pair table with the
# full 3x3 cross join across ECs/FAPs/sCs sending & target clusters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Tgfb2","Tgfbr3","ECs",2,0.6666666666666666,1.465615333333333,4.396846,0.04672291954663727,0.04672291954663728,3,1,14.85604233333333,44.568127,0.09286934904108346,0.09286934904108346,21.77324343638244,195.959190927442,0.004339127123595118,0.004339127123595118),
    @("ECs","Tgfb2","Tgfbr3","FAPs",2,0.6666666666666666,1.465615333333333,4.396846,0.04672291954663727,0.04672291954663728,3,1,114.2734143333333,342.820243,0.7143556381787382,0.7143556381787382,167.4808682392865,1507.327814153578,0.03337678101031191,0.03337678101031191),
    @("ECs","Tgfb2","Tgfbr3","sCs",2,0.6666666666666666,1.465615333333333,4.396846,0.04672291954663727,0.04672291954663728,3,1,30.83766366666667,92.512991,0.1927750127801784,0.1927750127801784,45.19615271404289,406.765374426386,0.009007011412730246,0.009007011412730248),
    @("FAPs","Tgfb2","Tgfbr3","ECs",3,1,18.88237266666667,56.64711800000001,0.6019584804341267,0.6019584804341268,3,1,14.85604233333333,44.568127,0.09286934904108346,0.09286934904108346,280.5173276897763,2524.655949207986,0.05590349222767712,0.05590349222767713),
    @("FAPs","Tgfb2","Tgfbr3","FAPs",3,1,18.88237266666667,56.64711800000001,0.6019584804341267,0.6019584804341268,3,1,114.2734143333333,342.820243,0.7143556381787382,0.7143556381787382,2157.753195334409,19419.77875800968,0.430012434447624,0.4300124344476241),
    @("FAPs","Tgfb2","Tgfbr3","sCs",3,1,18.88237266666667,56.64711800000001,0.6019584804341267,0.6019584804341268,3,1,30.83766366666667,92.512991,0.1927750127801784,0.1927750127801784,582.2882575233265,5240.594317709939,0.1160425537588255,0.1160425537588256),
    @("sCs","Tgfb2","Tgfbr3","ECs",3,1,11.020243,33.060729,0.351318600019236,0.351318600019236,3,1,14.85604233333333,44.568127,0.09286934904108346,0.09286934904108346,163.7171965316203,1473.454768784583,0.03262672968981122,0.03262672968981122),
    @("sCs","Tgfb2","Tgfbr3","FAPs",3,1,11.020243,33.060729,0.351318600019236,0.351318600019236,3,1,114.2734143333333,342.820243,0.7143556381787382,0.7143556381787382,1259.320794393017,11333.88714953715,0.2509664227208022,0.2509664227208022),
    @("sCs","Tgfb2","Tgfbr3","sCs",3,1,11.020243,33.060729,0.351318600019236,0.351318600019236,3,1,30.83766366666667,92.512991,0.1927750127801784,0.1927750127801784,339.8385471589377,3058.546924430439,0.0677254476086226,0.0677254476086226)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
